$wb = $excel.ActiveWorkbook

# --- Sheet1 (TC101): column widths ---
$ws1 = $wb.Worksheets.Item("TC101")
$ws1.Columns.Item(1).ColumnWidth = 15.499999999999998
$ws1.Columns.Item(2).ColumnWidth = 16.666666666666668

# --- Recreate TC102 so it gets a fresh internal sheetId (2 -> 3), same name/position ---
$temp = $wb.Worksheets.Add()
$wb.Worksheets.Item("TC102").Delete() | Out-Null
$temp.Name = "TC102"
$temp.Move(, $ws1)
$ws2 = $wb.Worksheets.Item("TC102")

# --- Populate TC102 with the hotel-booking table (write order matters for shared-string order) ---
$ws2.Range("D1").Value = "No.of rooms"
$ws2.Range("E1").Value = "check in date"
$ws2.Range("F1").Value = "Check out date"
$ws2.Range("G2").Value = "2 - two"
$ws2.Range("A1").Value = "LOCATION"
$ws2.Range("B1").Value = "HOTELS"
$ws2.Range("C1").Value = "ROOM_TYPE"
$ws2.Range("G1").Value = "ADULTS PER ROOM"
$ws2.Range("H1").Value = "CHILD PER ROOM"
$ws2.Range("A2").Value = "Sydney"
$ws2.Range("B2").Value = "Hotel Creek"
$ws2.Range("D2").Value = "1 - One"
$ws2.Range("C2").Value = "Standard"
$ws2.Range("H2").Value = "2 - two"

$ws2.Range("E2").Value = 45055
$ws2.Range("E2").NumberFormat = "mm-dd-yy"
$ws2.Range("E2").Copy() | Out-Null
$ws2.Range("F2").PasteSpecial(-4122) | Out-Null
$ws2.Range("F2").Value = 45053

# --- Column widths on TC102 ---
$ws2.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(2).ColumnWidth = 18.166666666666668
$ws2.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws2.Columns.Item(4).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws2.Columns.Item(6).ColumnWidth = 15.0
$ws2.Columns.Item(7).ColumnWidth = 14.333333333333332
$ws2.Columns.Item(8).ColumnWidth = 17.666666666666668

# --- Selections: TC101 -> L19 (no longer the active tab), TC102 -> G6 (active tab) ---
$ws1.Range("L19").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("G6").Select() | Out-Null
